$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$s.Shapes.Item(2).TextFrame.TextRange.Paragraphs(13).Text = "Most divisive movies between audience and top critic"
